$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.149.25"
$ws.Range("E2").Value = "  +3.24%  "
$ws.Range("D3").Value = "2.537.91"
$ws.Range("E3").Value = "  +4.87%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.996"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "527.16"
$ws.Range("E5").Value = "  +2.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.15"
$ws.Range("E6").Value = "  +5.03%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +3.52%  "
$ws.Range("D9").Value = "2.536.52"
$ws.Range("E9").Value = "  +4.59%  "
$ws.Range("E10").Value = "  +3.95%  "
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("D14").Value = "2.964.19"
$ws.Range("E14").Value = "  +3.96%  "
$ws.Range("D15").Value = "58.962.30"
$ws.Range("E15").Value = "  +2.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.38"
$ws.Range("E16").Value = "  +4.41%  "
$ws.Range("E17").Value = "  +3.87%  "
$ws.Range("D18").Value = "2.516.23"
$ws.Range("E18").Value = "  +3.76%  "
$ws.Range("E19").Value = "  +3.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "324.00"
$ws.Range("E20").Value = "  +3.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.20"
$ws.Range("E21").Value = "  +3.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.10"
$ws.Range("E22").Value = "  +8.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.996"
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.01"
$ws.Range("E24").Value = "  +2.39%  "
$ws.Range("E25").Value = "  +2.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +1.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.52"
$ws.Range("E28").Value = "  +4.79%  "
$ws.Range("D29").Value = "0.0₃0762"
$ws.Range("E29").Value = "  +6.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.24"
$ws.Range("E30").Value = "  +7.68%  "
$ws.Range("E31").Value = "  +4.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "169.72"
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("E33").Value = "  +3.19%  "
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.31"
$ws.Range("E36").Value = "  +3.73%  "
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.00"
$ws.Range("E38").Value = "  +3.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.51"
$ws.Range("E39").Value = "  +5.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.72"
$ws.Range("E40").Value = "  +1.23%  "
$ws.Range("E41").Value = "  +3.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "280.62"
$ws.Range("E42").Value = "  +6.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "134.87"
$ws.Range("E43").Value = "  +11.44%  "
$ws.Range("E44").Value = "  +4.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.10"
$ws.Range("E45").Value = "  +5.19%  "
$ws.Range("E47").Value = "  +2.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0507"
$ws.Range("E48").Value = "  +5.64%  "
$ws.Range("E49").Value = "  +4.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.18"
$ws.Range("E50").Value = "  +4.43%  "
$ws.Range("D51").Value = "1.757.91"
$ws.Range("E51").Value = "  +3.94%  "
